# Automatic update of files.
# Increment the "Förändrad" (Changed) date in column C for rows 2-5 by one day
# (from 2023-09-14 to 2023-09-15), matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2..5) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    $cell.Value2 = $cell.Value2 + 1
}
